$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2260,1).Value = '2025-07-25T11:01:23+00:00'
$ws.Cells.Item(2260,2).Value = 'EXTERNAL:- BBC Radio 5 Live - Wk32 - 2025-08-10 - Sunday'
$ws.Cells.Item(2260,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2260,4).Value = $true

$ws.Cells.Item(2261,1).Value = '2025-07-25T11:01:22+00:00'
$ws.Cells.Item(2261,2).Value = 'EXTERNAL:- BBC Three - Wk31 - 2025-08-06 - Wednesday'
$ws.Cells.Item(2261,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2261,4).Value = $true

$ws.Cells.Item(2262,1).Value = '2025-07-25T11:01:22+00:00'
$ws.Cells.Item(2262,2).Value = 'EXTERNAL:- BBC Two HD - Wk31 - 2025-08-05 - Tuesday'
$ws.Cells.Item(2262,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2262,4).Value = $true

$ws.Cells.Item(2263,1).Value = '2025-07-25T11:01:21+00:00'
$ws.Cells.Item(2263,2).Value = 'EXTERNAL:- BBC Two HD - Wk31 - 2025-08-05 - Tuesday'
$ws.Cells.Item(2263,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2263,4).Value = $true

$ws.Cells.Item(2264,1).Value = '2025-07-25T11:01:12+00:00'
$ws.Cells.Item(2264,2).Value = 'EXTERNAL:- BBC Asian Network - Wk32 - 2025-08-09 - Saturday'
$ws.Cells.Item(2264,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2264,4).Value = $true

$ws.Cells.Item(2265,1).Value = '2025-07-25T11:01:12+00:00'
$ws.Cells.Item(2265,2).Value = 'EXTERNAL:- BBC Radio 3 - Wk32 - 2025-08-13 - Wednesday'
$ws.Cells.Item(2265,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2265,4).Value = $true

$ws.Cells.Item(2266,1).Value = '2025-07-25T11:01:10+00:00'
$ws.Cells.Item(2266,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk32 - 2025-08-13 - Wednesday'
$ws.Cells.Item(2266,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2266,4).Value = $true

$ws.Cells.Item(2267,1).Value = '2025-07-25T11:01:09+00:00'
$ws.Cells.Item(2267,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk31 - 2025-08-08 - Friday'
$ws.Cells.Item(2267,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2267,4).Value = $true

$ws.Cells.Item(2268,1).Value = '2025-07-25T11:01:08+00:00'
$ws.Cells.Item(2268,2).Value = 'EXTERNAL:- BBC Radio nan Gàidheal - Wk32 - 2025-08-13 - Wednesday'
$ws.Cells.Item(2268,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2268,4).Value = $true

$ws.Cells.Item(2269,1).Value = '2025-07-25T11:01:02+00:00'
$ws.Cells.Item(2269,2).Value = 'EXTERNAL:- BBC Two HD - Wk31 - 2025-08-08 - Friday'
$ws.Cells.Item(2269,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2269,4).Value = $true

$ws.Cells.Item(2270,1).Value = '2025-07-25T11:01:02+00:00'
$ws.Cells.Item(2270,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk32 - 2025-08-10 - Sunday'
$ws.Cells.Item(2270,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2270,4).Value = $true

$ws.Cells.Item(2271,1).Value = '2025-07-25T11:01:01+00:00'
$ws.Cells.Item(2271,2).Value = 'EXTERNAL:- BBC Two HD - Wk31 - 2025-08-08 - Friday'
$ws.Cells.Item(2271,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2271,4).Value = $true

$ws.Cells.Item(2272,1).Value = '2025-07-25T11:01:00+00:00'
$ws.Cells.Item(2272,2).Value = 'EXTERNAL:- BBC Radio 4 FM - Wk31 - 2025-08-08 - Friday'
$ws.Cells.Item(2272,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2272,4).Value = $true

$ws.Cells.Item(2273,1).Value = '2025-07-25T11:00:56+00:00'
$ws.Cells.Item(2273,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk31 - 2025-08-07 - Thursday'
$ws.Cells.Item(2273,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2273,4).Value = $true

$ws.Cells.Item(2274,1).Value = '2025-07-25T11:00:56+00:00'
$ws.Cells.Item(2274,2).Value = 'EXTERNAL:- BBC Radio 3 - Wk32 - 2025-08-12 - Tuesday'
$ws.Cells.Item(2274,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2274,4).Value = $true

$ws.Cells.Item(2275,1).Value = '2025-07-25T11:00:55+00:00'
$ws.Cells.Item(2275,2).Value = 'EXTERNAL:- BBC Radio 4 Extra - Wk31 - 2025-08-05 - Tuesday'
$ws.Cells.Item(2275,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2275,4).Value = $true

$ws.Cells.Item(2276,1).Value = '2025-07-25T11:00:53+00:00'
$ws.Cells.Item(2276,2).Value = 'EXTERNAL:- BBC Radio 4 FM - Wk31 - 2025-08-07 - Thursday'
$ws.Cells.Item(2276,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2276,4).Value = $true

$ws.Cells.Item(2277,1).Value = '2025-07-25T11:00:52+00:00'
$ws.Cells.Item(2277,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk32 - 2025-08-12 - Tuesday'
$ws.Cells.Item(2277,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2277,4).Value = $true

$ws.Cells.Item(2278,1).Value = '2025-07-25T11:00:52+00:00'
$ws.Cells.Item(2278,2).Value = 'EXTERNAL:- BBC Two HD - Wk31 - 2025-08-06 - Wednesday'
$ws.Cells.Item(2278,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2278,4).Value = $true

$ws.Cells.Item(2279,1).Value = '2025-07-25T11:00:52+00:00'
$ws.Cells.Item(2279,2).Value = 'EXTERNAL:- BBC Radio 3 - Wk32 - 2025-08-14 - Thursday'
$ws.Cells.Item(2279,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2279,4).Value = $true

$ws.Cells.Item(2280,1).Value = '2025-07-25T11:00:51+00:00'
$ws.Cells.Item(2280,2).Value = 'EXTERNAL:- BBC Two HD - Wk31 - 2025-08-06 - Wednesday'
$ws.Cells.Item(2280,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2280,4).Value = $true

$ws.Cells.Item(2281,1).Value = '2025-07-25T11:00:48+00:00'
$ws.Cells.Item(2281,2).Value = 'EXTERNAL:- BBC Three - Wk31 - 2025-08-07 - Thursday'
$ws.Cells.Item(2281,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2281,4).Value = $true

$ws.Cells.Item(2282,1).Value = '2025-07-25T11:00:46+00:00'
$ws.Cells.Item(2282,2).Value = 'EXTERNAL:- BBC Radio 4 FM - Wk32 - 2025-08-13 - Wednesday'
$ws.Cells.Item(2282,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2282,4).Value = $true

$ws.Cells.Item(2283,1).Value = '2025-07-25T11:00:46+00:00'
$ws.Cells.Item(2283,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk31 - 2025-08-06 - Wednesday'
$ws.Cells.Item(2283,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2283,4).Value = $true

$ws.Cells.Item(2284,1).Value = '2025-07-25T11:00:45+00:00'
$ws.Cells.Item(2284,2).Value = 'EXTERNAL:- BBC Radio 1 - Wk32 - 2025-08-10 - Sunday'
$ws.Cells.Item(2284,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2284,4).Value = $true

$ws.Cells.Item(2285,1).Value = '2025-07-25T11:00:45+00:00'
$ws.Cells.Item(2285,2).Value = 'EXTERNAL:- BBC Radio 4 FM - Wk32 - 2025-08-10 - Sunday'
$ws.Cells.Item(2285,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2285,4).Value = $true

$ws.Cells.Item(2286,1).Value = '2025-07-25T11:00:43+00:00'
$ws.Cells.Item(2286,2).Value = 'EXTERNAL:- BBC Radio 4 FM - Wk32 - 2025-08-14 - Thursday'
$ws.Cells.Item(2286,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2286,4).Value = $true

$ws.Cells.Item(2287,1).Value = '2025-07-25T11:00:42+00:00'
$ws.Cells.Item(2287,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk32 - 2025-08-11 - Monday'
$ws.Cells.Item(2287,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2287,4).Value = $true

$ws.Cells.Item(2288,1).Value = '2025-07-25T11:00:42+00:00'
$ws.Cells.Item(2288,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk32 - 2025-08-09 - Saturday'
$ws.Cells.Item(2288,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2288,4).Value = $true

$ws.Cells.Item(2289,1).Value = '2025-07-25T11:00:42+00:00'
$ws.Cells.Item(2289,2).Value = 'EXTERNAL:- BBC Radio 3 - Wk32 - 2025-08-11 - Monday'
$ws.Cells.Item(2289,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2289,4).Value = $true

$ws.Cells.Item(2290,1).Value = '2025-07-25T11:00:38+00:00'
$ws.Cells.Item(2290,2).Value = 'EXTERNAL:- BBC Radio 5 Live - Wk31 - 2025-08-05 - Tuesday'
$ws.Cells.Item(2290,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2290,4).Value = $true

$ws.Cells.Item(2291,1).Value = '2025-07-25T11:00:37+00:00'
$ws.Cells.Item(2291,2).Value = 'EXTERNAL:- BBC Radio 4 FM - Wk32 - 2025-08-12 - Tuesday'
$ws.Cells.Item(2291,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2291,4).Value = $true

$ws.Cells.Item(2292,1).Value = '2025-07-25T11:00:36+00:00'
$ws.Cells.Item(2292,2).Value = 'EXTERNAL:- BBC Radio 5 Sports Extra - Wk31 - 2025-08-05 - Tuesday'
$ws.Cells.Item(2292,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2292,4).Value = $true

$ws.Cells.Item(2293,1).Value = '2025-07-25T11:00:35+00:00'
$ws.Cells.Item(2293,2).Value = 'EXTERNAL:- BBC Radio 3 - Wk32 - 2025-08-10 - Sunday'
$ws.Cells.Item(2293,3).Value = 'pressportal@bbc.co.uk'
$ws.Cells.Item(2293,4).Value = $true

$ws.Cells.Item(2294,1).Value = '2025-07-25T11:00:34+00:00'
$ws.Cells.Item(2294,2).Value = 'EXTERNAL:- More4 week 33 provisional listings available'
$ws.Cells.Item(2294,3).Value = 'FHowarth@Channel4.co.uk'
$ws.Cells.Item(2294,4).Value = $false

$lo = $ws.ListObjects.Item(1)
$newRange = $ws.Range("A1:D2294")
$lo.Resize($newRange)

Write-Host "Done. Table range:" $lo.Range.Address()
Write-Host "Dimension check A2294:" $ws.Cells.Item(2294,1).Value()